$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/6/2025  Through  1/12/2025"

# --- Data table updates ---
# Cells whose style does not change: direct value assignment
$ws.Range("I15").Value = 2
$ws.Range("K15").Value = 100
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 4
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 64.285714285714
$ws.Range("I16").Value = 11
$ws.Range("L16").Value = 83.333333333333
$ws.Range("M16").Value = 22.222222222222
$ws.Range("N16").Value = 10
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -13.636363636363
$ws.Range("I17").Value = 6
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = -33.333333333333
$ws.Range("L17").Value = -45.454545454545
$ws.Range("M17").Value = -14.285714285714
$ws.Range("N17").Value = -25
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = -50
$ws.Range("N18").Value = -90.909090909090
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 37.5
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -37.254901960784
$ws.Range("I19").Value = 15
$ws.Range("J19").Value = 19
$ws.Range("K19").Value = -21.052631578947
$ws.Range("L19").Value = -25
$ws.Range("M19").Value = 15.384615384615
$ws.Range("N19").Value = 7.142857142857
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 37.5
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = -8.571428571428
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 31.25
$ws.Range("L20").Value = 75
$ws.Range("N20").Value = -77.173913043478
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 19.230769230769
$ws.Range("F21").Value = 119
$ws.Range("G21").Value = 133
$ws.Range("H21").Value = -10.526315789473
$ws.Range("I21").Value = 58
$ws.Range("J21").Value = 57
$ws.Range("K21").Value = 1.754385964912
$ws.Range("L21").Value = 16
$ws.Range("M21").Value = 41.463414634146
$ws.Range("N21").Value = -63.291139240506
$ws.Range("M22").Value = 0
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = -50
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 76.470588235294
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = -28.037383177570
$ws.Range("I24").Value = 38
$ws.Range("J24").Value = 37
$ws.Range("K24").Value = 2.702702702702
$ws.Range("L24").Value = -24
$ws.Range("M24").Value = 15.151515151515
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = -52.112676056338
$ws.Range("I25").Value = 15
$ws.Range("J25").Value = 25
$ws.Range("K25").Value = -40
$ws.Range("L25").Value = -58.333333333333
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -21.052631578947
$ws.Range("I26").Value = 17
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = 21.428571428571
$ws.Range("L26").Value = 21.428571428571
$ws.Range("M26").Value = 13.333333333333
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 2
$ws.Range("K27").Value = 100
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = 500
$ws.Range("H29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = -50
$ws.Range("J43").Value = 179
$ws.Range("K43").Value = -51.621621621621
$ws.Range("L43").Value = -53.746770025839
$ws.Range("M43").Value = -77.708592777085
$ws.Range("N43").Value = -78.276699029126
$ws.Range("J44").Value = 683
$ws.Range("K44").Value = 44.703389830508
$ws.Range("L44").Value = 65.776699029126
$ws.Range("M44").Value = 64.578313253012
$ws.Range("N44").Value = 44.397463002114
$ws.Range("J46").Value = 1862
$ws.Range("K46").Value = -1.637612255678
$ws.Range("L46").Value = -24.828421477593
$ws.Range("M46").Value = -54.915254237288
$ws.Range("N46").Value = -60.609265919187

# Cells becoming numeric (from text "N/A" style 13) -> set value, then fix style via PasteSpecial(formats)
$ws.Range("M15").Value = 100
$ws.Range("H14").Copy()
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("D16").Value = 5
$ws.Range("G14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = -20
$ws.Range("H14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("J16").Value = 5
$ws.Range("G14").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$ws.Range("K16").Value = 120
$ws.Range("H14").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$ws.Range("L18").Value = 200
$ws.Range("H14").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("M18").Value = -57.142857142857
$ws.Range("H14").Copy()
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M20").Value = 425
$ws.Range("H14").Copy()
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("I22").Value = 1
$ws.Range("G14").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("L28").Value = 200
$ws.Range("H14").Copy()
$ws.Range("L28").PasteSpecial(-4122)
$ws.Range("C29").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1
$ws.Range("G14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("I29").Value = 1
$ws.Range("G14").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("C30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("I30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("I30").PasteSpecial(-4122)
$ws.Range("J33").Value = 1
$ws.Range("G14").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("K33").Value = -100
$ws.Range("H14").Copy()
$ws.Range("K33").PasteSpecial(-4122)

# Cells becoming text "N/A" (style 13) -> set text value (forced via leading apostrophe), then fix style via PasteSpecial(formats)
$ws.Range("D14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
